$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 858.5
$ws.Cells.Item(18, 9).Value = 960.4
$ws.Cells.Item(18, 10).Value = 349
$ws.Cells.Item(18, 11).Value = 960.4
$ws.Cells.Item(18, 12).Value = 349
$ws.Cells.Item(18, 13).Value = -676.4
$ws.Cells.Item(18, 14).Value = -917
$ws.Cells.Item(28, 8).Value = 939
$ws.Cells.Item(28, 9).Value = 514.6667
$ws.Cells.Item(28, 11).Value = 514.6667
$ws.Cells.Item(28, 13).Value = -29.66669999999999
$ws.Cells.Item(32, 8).Value = 3999
$ws.Cells.Item(32, 10).Value = 3999
$ws.Cells.Item(32, 12).Value = 3999
$ws.Cells.Item(32, 14).Value = -4651
$ws.Cells.Item(33, 8).Value = 73045
$ws.Cells.Item(33, 9).Value = 85135
$ws.Cells.Item(33, 11).Value = 85135
$ws.Cells.Item(33, 13).Value = -84906
$ws.Cells.Item(42, 8).Value = 646.5
$ws.Cells.Item(42, 9).Value = 391.6
$ws.Cells.Item(42, 10).Value = 828.5714
$ws.Cells.Item(42, 11).Value = 1174.8
$ws.Cells.Item(42, 12).Value = 2485.7142
$ws.Cells.Item(42, 13).Value = -944.8000000000002
$ws.Cells.Item(42, 14).Value = -2945.7142
$ws.Cells.Item(43, 8).Value = 4544.3193
$ws.Cells.Item(43, 9).Value = 4934.8
$ws.Cells.Item(43, 10).Value = 4100.591
$ws.Cells.Item(43, 11).Value = 4934.8
$ws.Cells.Item(43, 12).Value = 4100.591
$ws.Cells.Item(43, 13).Value = -4865.8
$ws.Cells.Item(43, 14).Value = -4238.591
$ws.Cells.Item(62, 8).Value = 4387.3
$ws.Cells.Item(62, 9).Value = 4398.625
$ws.Cells.Item(62, 10).Value = 4342
$ws.Cells.Item(62, 11).Value = 4398.625
$ws.Cells.Item(62, 12).Value = 4342
$ws.Cells.Item(62, 13).Value = -3774.625
$ws.Cells.Item(62, 14).Value = -5590
$ws.Cells.Item(64, 8).Value = 5951.2
$ws.Cells.Item(64, 9).Value = 6778.5
$ws.Cells.Item(64, 10).Value = 5399.6665
$ws.Cells.Item(64, 11).Value = 6778.5
$ws.Cells.Item(64, 12).Value = 5399.6665
$ws.Cells.Item(64, 13).Value = -6530.5
$ws.Cells.Item(64, 14).Value = -5895.6665
$ws.Cells.Item(65, 8).Value = 4387.3
$ws.Cells.Item(65, 9).Value = 4398.625
$ws.Cells.Item(65, 10).Value = 4342
$ws.Cells.Item(65, 11).Value = 21993.125
$ws.Cells.Item(65, 12).Value = 21710
$ws.Cells.Item(65, 13).Value = -18873.125
$ws.Cells.Item(65, 14).Value = -27950
$ws.Cells.Item(67, 8).Value = 5951.2
$ws.Cells.Item(67, 9).Value = 6778.5
$ws.Cells.Item(67, 10).Value = 5399.6665
$ws.Cells.Item(67, 11).Value = 6778.5
$ws.Cells.Item(67, 12).Value = 5399.6665
$ws.Cells.Item(67, 13).Value = -5920.5
$ws.Cells.Item(67, 14).Value = -7115.6665
$ws.Cells.Item(86, 8).Value = 2451.6
$ws.Cells.Item(86, 9).Value = 3117.3333
$ws.Cells.Item(86, 10).Value = 1453
$ws.Cells.Item(86, 11).Value = 3117.3333
$ws.Cells.Item(86, 12).Value = 1453
$ws.Cells.Item(86, 13).Value = -1994.3333
$ws.Cells.Item(86, 14).Value = -3699
$ws.Cells.Item(88, 8).Value = 41400
$ws.Cells.Item(88, 10).Value = 41400
$ws.Cells.Item(88, 12).Value = 41400
$ws.Cells.Item(88, 14).Value = -42212
$ws.Cells.Item(89, 8).Value = 2451.6
$ws.Cells.Item(89, 9).Value = 3117.3333
$ws.Cells.Item(89, 10).Value = 1453
$ws.Cells.Item(89, 11).Value = 15586.6665
$ws.Cells.Item(89, 12).Value = 7265
$ws.Cells.Item(89, 13).Value = -9970.666499999999
$ws.Cells.Item(89, 14).Value = -18497
$ws.Cells.Item(91, 8).Value = 41400
$ws.Cells.Item(91, 10).Value = 41400
$ws.Cells.Item(91, 12).Value = 41400
$ws.Cells.Item(91, 14).Value = -44208
$ws.Cells.Item(111, 8).Value = 1665.5714
$ws.Cells.Item(111, 9).Value = 1631.8
$ws.Cells.Item(111, 11).Value = 4895.4
$ws.Cells.Item(111, 13).Value = -1828.4
$ws.Cells.Item(141, 8).Value = 3718.1
$ws.Cells.Item(141, 9).Value = 3718.1
$ws.Cells.Item(141, 11).Value = 11154.3
$ws.Cells.Item(141, 13).Value = -5974.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 1999.3334
$ws.Cells.Item(88, 9).Value = 1932.6666
$ws.Cells.Item(88, 11).Value = 1932.6666
$ws.Cells.Item(88, 13).Value = -1526.6666
$ws.Cells.Item(91, 8).Value = 1999.3334
$ws.Cells.Item(91, 9).Value = 1932.6666
$ws.Cells.Item(91, 11).Value = 1932.6666
$ws.Cells.Item(91, 13).Value = -528.6666
$ws.Cells.Item(132, 8).Value = 33338538
$ws.Cells.Item(132, 9).Value = 4466.077
$ws.Cells.Item(132, 11).Value = 13398.231
$ws.Cells.Item(132, 13).Value = -10868.231
$ws.Cells.Item(139, 8).Value = 58999
$ws.Cells.Item(139, 10).Value = 58999
$ws.Cells.Item(139, 12).Value = 58999
$ws.Cells.Item(139, 14).Value = -69279

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3740.8096
$ws.Cells.Item(20, 9).Value = 3321.7144
$ws.Cells.Item(20, 10).Value = 3950.3572
$ws.Cells.Item(20, 11).Value = 3321.7144
$ws.Cells.Item(20, 12).Value = 3950.3572
$ws.Cells.Item(20, 13).Value = -3074.7144
$ws.Cells.Item(20, 14).Value = -4444.3572
$ws.Cells.Item(39, 8).Value = 21220.2
$ws.Cells.Item(39, 10).Value = 21220.2
$ws.Cells.Item(39, 12).Value = 21220.2
$ws.Cells.Item(39, 14).Value = -21998.2
$ws.Cells.Item(54, 8).Value = 7881.25
$ws.Cells.Item(54, 9).Value = 7881.25
$ws.Cells.Item(54, 11).Value = 7881.25
$ws.Cells.Item(54, 13).Value = -7397.25
$ws.Cells.Item(107, 8).Value = 15849.25
$ws.Cells.Item(107, 9).Value = 17970.715
$ws.Cells.Item(107, 10).Value = 999
$ws.Cells.Item(107, 11).Value = 17970.715
$ws.Cells.Item(107, 12).Value = 999
$ws.Cells.Item(107, 13).Value = -16050.715
$ws.Cells.Item(107, 14).Value = -4839

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(39, 8).Value = 52250.5
$ws.Cells.Item(39, 10).Value = 52250.5
$ws.Cells.Item(39, 12).Value = 52250.5
$ws.Cells.Item(39, 14).Value = -53032.5
$ws.Cells.Item(49, 8).Value = 52250.5
$ws.Cells.Item(49, 10).Value = 52250.5
$ws.Cells.Item(49, 12).Value = 52250.5
$ws.Cells.Item(49, 14).Value = -52614.5
$ws.Cells.Item(122, 8).Value = 44266.266
$ws.Cells.Item(122, 9).Value = 6124.75
$ws.Cells.Item(122, 10).Value = 87856.57000000001
$ws.Cells.Item(122, 11).Value = 18374.25
$ws.Cells.Item(122, 12).Value = 263569.71
$ws.Cells.Item(122, 13).Value = -15924.25
$ws.Cells.Item(122, 14).Value = -268469.71
$ws.Cells.Item(132, 8).Value = 4229.273
$ws.Cells.Item(132, 9).Value = 5308
$ws.Cells.Item(132, 10).Value = 2934.8
$ws.Cells.Item(132, 11).Value = 15924
$ws.Cells.Item(132, 12).Value = 8804.400000000001
$ws.Cells.Item(132, 13).Value = -13394
$ws.Cells.Item(132, 14).Value = -13864.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 72.29412000000001
$ws.Cells.Item(17, 9).Value = 68
$ws.Cells.Item(17, 11).Value = 204
$ws.Cells.Item(17, 13).Value = -35
$ws.Cells.Item(34, 8).Value = 4627.478
$ws.Cells.Item(34, 10).Value = 6788.533
$ws.Cells.Item(34, 12).Value = 20365.599
$ws.Cells.Item(34, 14).Value = -20533.599
$ws.Cells.Item(39, 8).Value = 6213.5713
$ws.Cells.Item(39, 10).Value = 7780.8
$ws.Cells.Item(39, 12).Value = 23342.4
$ws.Cells.Item(39, 14).Value = -23930.4
$ws.Cells.Item(55, 8).Value = 4699
$ws.Cells.Item(55, 10).Value = 6999.1665
$ws.Cells.Item(55, 12).Value = 20997.4995
$ws.Cells.Item(55, 14).Value = -21351.4995
$ws.Cells.Item(112, 8).Value = 6545.4546
$ws.Cells.Item(118, 8).Value = 2506.1428
$ws.Cells.Item(118, 9).Value = 2590.5
$ws.Cells.Item(118, 11).Value = 7771.5
$ws.Cells.Item(118, 13).Value = -6528.5
$ws.Cells.Item(140, 8).Value = 1989.381
$ws.Cells.Item(140, 9).Value = 1741.35
$ws.Cells.Item(140, 11).Value = 5224.049999999999
$ws.Cells.Item(140, 13).Value = -44.04999999999927

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 4750
$ws.Cells.Item(132, 9).Value = 4750
$ws.Cells.Item(132, 11).Value = 14250
$ws.Cells.Item(132, 13).Value = -11720

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 5547.9165
$ws.Cells.Item(22, 9).Value = 900
$ws.Cells.Item(22, 10).Value = 6477.5
$ws.Cells.Item(22, 11).Value = 900
$ws.Cells.Item(22, 12).Value = 6477.5
$ws.Cells.Item(22, 13).Value = -605
$ws.Cells.Item(22, 14).Value = -7067.5
$ws.Cells.Item(27, 8).Value = 5547.9165
$ws.Cells.Item(27, 9).Value = 900
$ws.Cells.Item(27, 10).Value = 6477.5
$ws.Cells.Item(27, 11).Value = 900
$ws.Cells.Item(27, 12).Value = 6477.5
$ws.Cells.Item(27, 13).Value = -793
$ws.Cells.Item(27, 14).Value = -6691.5
$ws.Cells.Item(39, 8).Value = 47532.5
$ws.Cells.Item(39, 10).Value = 70065
$ws.Cells.Item(39, 12).Value = 70065
$ws.Cells.Item(39, 14).Value = -70985
$ws.Cells.Item(122, 8).Value = 7817.35
$ws.Cells.Item(122, 9).Value = 7518.467
$ws.Cells.Item(122, 10).Value = 8714
$ws.Cells.Item(122, 11).Value = 22555.401
$ws.Cells.Item(122, 12).Value = 26142
$ws.Cells.Item(122, 13).Value = -20105.401
$ws.Cells.Item(122, 14).Value = -31042

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3666.6667
$ws.Cells.Item(62, 10).Value = 3666.6667
$ws.Cells.Item(62, 12).Value = 3666.6667
$ws.Cells.Item(62, 14).Value = -4914.6667
$ws.Cells.Item(65, 8).Value = 3666.6667
$ws.Cells.Item(65, 10).Value = 3666.6667
$ws.Cells.Item(65, 12).Value = 18333.3335
$ws.Cells.Item(65, 14).Value = -24573.3335
$ws.Cells.Item(96, 8).Value = 2666.5833
$ws.Cells.Item(96, 10).Value = 2100
$ws.Cells.Item(96, 12).Value = 2100
$ws.Cells.Item(96, 14).Value = -4846
$ws.Cells.Item(107, 8).Value = 1939.8214
$ws.Cells.Item(107, 9).Value = 1466.1
$ws.Cells.Item(107, 10).Value = 3124.125
$ws.Cells.Item(107, 11).Value = 4398.299999999999
$ws.Cells.Item(107, 12).Value = 9372.375
$ws.Cells.Item(107, 13).Value = -2478.299999999999
$ws.Cells.Item(107, 14).Value = -13212.375
$ws.Cells.Item(122, 8).Value = 2367.25
$ws.Cells.Item(122, 9).Value = 2234.5
$ws.Cells.Item(122, 11).Value = 6703.5
$ws.Cells.Item(122, 13).Value = -4253.5
$ws.Cells.Item(126, 8).Value = 3372.6316
$ws.Cells.Item(126, 9).Value = 3322.7
$ws.Cells.Item(126, 11).Value = 9968.099999999999
$ws.Cells.Item(126, 13).Value = -7498.099999999999
